$wb = $excel.ActiveWorkbook
$about = $wb.Worksheets.Item(1)
$mcf = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 0) "MCF" sheet: update the forced outage rate assumptions for hard coal and
#    nuclear per the cited NRDC source. Done first so the final active sheet
#    / selection ends up on "About", matching the saved workbook state.
# ---------------------------------------------------------------------------
$mcf.Range("B2").Value = 0.85
$mcf.Range("B5").Value = 0.91
$mcf.Activate()
$mcf.Range("B4").Select()

# ---------------------------------------------------------------------------
# 1) "About" sheet: rework the Source block into a full citation, and extend
#    the Notes block with the new forced-outage-rate explanation.
# ---------------------------------------------------------------------------
$about.Activate()

# Make room for the expanded source/citation block: rows 4-9 (old) need to
# become rows 7-12 (new), so insert 3 fresh rows above the old row 4.
$about.Rows("4:6").Insert()

# --- Row 3: "Source:" label stays, but the value becomes the article title.
$about.Rows(3).ClearFormats()
$about.Range("A3").Font.Bold = $true
$about.Range("B3").Value = "The Myth of the 24/7/365 Power Plant"
$about.Range("E3").Font.Bold = $true
$about.Range("H3").Font.Bold = $true

# --- Row 4: publication year.
$about.Range("A4").Font.Bold = $true
$about.Range("B4").Value = 2019
$about.Range("B4").HorizontalAlignment = -4131
$about.Range("E4").Font.Bold = $true
$about.Range("H4").Font.Bold = $true

# --- Row 5: publisher / author organization.
$about.Range("A5").Font.Bold = $true
$about.Range("B5").Font.Bold = $false
$about.Range("B5").Interior.Pattern = 0
$about.Range("B5").Value = "National Resources Defense Council"
$about.Range("E5").Font.Bold = $true
$about.Range("H5").Font.Bold = $true

# --- Row 6: the actual hyperlink to the source.
$about.Range("A6").Font.Bold = $true
$about.Hyperlinks.Add($about.Range("B6"), "https://www.nrdc.org/bio/rachel-fakhry/myth-247365-power-plant")
$about.Range("B6").Style = "Hyperlink"
$about.Range("E6").Font.Bold = $true
$about.Range("H6").Font.Bold = $true

# --- Row 11 (old row 8, shifted by the insert above): drop the leftover
#     uniform-row formatting so "a 5% penalty..." looks like the other note
#     lines, matching the simplified layout.
$about.Rows(11).ClearFormats()
$about.Range("H11").HorizontalAlignment = -4131

# --- Notes text changes: the old single closing sentence is replaced by two
#     sentences that cite NRDC and add the coal penalty explanation.
$about.Range("A12").Value = "year due to factors like plant maintenance, based on the value for gas plants"
$about.Range("A13").Value = "in the NRDC source above. For coal, we apply a 15% penalty based on NRDC."

# --- Final selection / dimension bookkeeping to match the edited layout.
$about.Range("A14").Select()
